$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Update row 8 values in place (date + quantity)
$ws1.Cells.Item(8, 1).Value = 45319.99999999999
$ws1.Cells.Item(8, 2).Value = 24

# Remove old rows 9-15 (now obsolete trailing data)
$ws1.Range("A9:B15").EntireRow.Delete()

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Update row 6 values in place (date + quantity)
$ws2.Cells.Item(6, 1).Value = 45322.99999999999
$ws2.Cells.Item(6, 2).Value = 24

# Remove old rows 7-9 (now obsolete trailing data)
$ws2.Range("A7:B9").EntireRow.Delete()
